$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level view tweak (window geometry) ---
# xWindow/windowWidth are cosmetic Excel-application window coordinates that
# are not exposed through this COM surface; skipped intentionally.

# --- Clear stale header/summary/date cells from the old 3-row dataset ---
$ws.Range("D1").ClearContents()
$ws.Range("E1").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

# --- Write the scraped LinkedIn job rows ---
# Row 1
$ws.Range("A1").Value = 'Kotlin Developer – Full-time, Remote'
$ws.Range("B1").Value = 'Toptal'
$ws.Range("C1").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F1").Value = 'https://www.linkedin.com/jobs/view/2648921423/?eBP=CwEAAAF6sFnSRCwq5yeVOz3lshQnKtwaXvffo3zRe0WIOMsL3YG-ZpZVIBMvLn9_9Gg9TCdDSdJZV0Hyw5FDjEMzMI5cz_FoWfwBkpuYjqAMaMSx3WQ6WkNg_dIky8ujmtX1wq7erWFq7TZcwDzQ3f6I4Xth_TqroVAQjD9seATVb1vo_xrPRjW9FPlzcfIYud-fUKvefkNkxN7DU_EowIKOo2LHfHedDFiMNeMfGwDQCBDPmA_6tDkp1__vfz1TgC2m-zE8hLBezIUMGVB5rzLWyc1fQ20UomzN7a_vfR0kYSzxF2zymghJYoCK2qf5kxMfru64SMbS7kubsyu_-OzkVzpysCRSGFbxAbYUQoSV_EivSouNmAFnljvdG83H&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=IQ08OXcRjo1ksy7zwvcxwA%3D%3D&trackingId=qmfYMMQ33yhTaibFpRih7g%3D%3D&trk=flagship3_search_srp_jobs'

# Row 2
$ws.Range("A2").Value = 'Senior C++ Developer – Full-time, Remote'
$ws.Range("B2").Value = 'Toptal'
$ws.Range("C2").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F2").Value = 'https://www.linkedin.com/jobs/view/2637116275/?eBP=CwEAAAF6sFnSRY3SupZg-UcGKeEKvkm77turW9ynZpV6gFC99HQkVQiCICnjCEFctuG0wopCyPPIRbdSJhvFdPUsbe3vXZbhRRbs92iChk4pMN1A6SHDFDRCpT46Ivzdtpu-G7Aa4PbC9moCmZs3N1AXSEkBL50ZcHrl5cjdzyMKiNIfqVWUiPya-2GaPmT3xvn6AOzknisDIwt4TN3UZTXIutPXsn_qfv2n28LVYKyZsOYh3JwQm5j-bk_KZAtVMr08zdrZKF0dkemg8vxet2utIjZn4k9MLJj8zuBpB_IgiV6RvnkBsvcSAMzbwxUzrFMRr-eq7cFIVWvlw_XtegcaTbXAs2h-RWS8E9Qg-IqGV8n46xPXgS9JBLBdfFyk&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=IQ08OXcRjo1ksy7zwvcxwA%3D%3D&trackingId=eRMZ5rfytaY3q5a4bfIUGQ%3D%3D&trk=flagship3_search_srp_jobs'

# Row 3
$ws.Range("A3").Value = 'Senior Software Engineer - iOS'
$ws.Range("B3").Value = 'HungerStation | هنقرستيشن'
$ws.Range("C3").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F3").Value = 'https://www.linkedin.com/jobs/view/2631668884/?eBP=CwEAAAF6sFnSRX4kM6Cb5n08ZX2qvHjeLgslll6OaPb5lD2xPJ9WjXUokMr6NSrfqO0vY2UTaffT1ii9SBAPXd25l5hS7XX7U0UVgiD2IX4PwUk_mUwNWAcY3GyIpBAkFS79XtjZls8ASZc0EEZsc4QKDocpWcv7JqNAwkOXg52ma1wQIDmpbYJdy66tvHv1urrhrix68GBg9787tVHuO7KCLYd2Chy4H93k8Rs5eGk5ug0Sf_FbzxBNwdQabZdUfnf5_1mr0N2ep0Lc4NWd5GIHZwyvSaUPt-MqTyk7Wx9a4HVN2bP80xnQseuop-mDyMS_FNG7MCzT8FkTmPYVshyZyIvghEcIZjVoyCvCMtO3D0ia5fVXxrvjODT4ur3Q&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=IQ08OXcRjo1ksy7zwvcxwA%3D%3D&trackingId=huPKmoZECG9Euxzb2aJ8zQ%3D%3D&trk=flagship3_search_srp_jobs'

# Row 4
$ws.Range("A4").Value = 'MuleSoft Developer (On-site)'
$ws.Range("B4").Value = 'Halian'
$ws.Range("C4").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("E4").Value = '
              2 days ago
 '
$ws.Range("F4").Value = 'https://www.linkedin.com/jobs/view/2643770873/?eBP=JOB_SEARCH_ORGANIC&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=IQ08OXcRjo1ksy7zwvcxwA%3D%3D&trackingId=petFiK%2BrK7UyfAmDkzxoFQ%3D%3D&trk=flagship3_search_srp_jobs'

# Row 5
$ws.Range("A5").Value = 'Senior Software Engineer I - Android'
$ws.Range("B5").Value = 'talabat'
$ws.Range("C5").Value = 'Business Bay, Dubai, United Arab Emirates'
$ws.Range("F5").Value = 'https://www.linkedin.com/jobs/view/2636396773/?eBP=CwEAAAF6sFnSRYP6wJETgDUt_BTRuRrvqb3UmnGfNJR10NZ5crRzyjrz53oHiu-VI9EfEilJ6BYTvbDgRL0iXWb0PJDn9gByXcRErBIOQyHfGWbTaahvz8hhYpCIoPARe6tqm9Xxl6Nm2JQ7YGiA_La1ZBxuGzCUdAZRBWjGK8X23N7Dyw9PQjGUc5XpJfhNXBE2NhB6JNEORIawBArmxrjAZKh3dBfgto7wM4HGNMQgd2O4_UYxdhh-qd83AmTAh6j9o4iefJIaFxexDF9B4sklfLjhjUYUE1Ze4w1xAdngzdAh4pvcrj7MzKpYCrX6Az0mYx3qt7vcJU4wT72iDGQnfAEeErvtiuQTubBJvp-5RDGqE5RVVFERZUj-og&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=IQ08OXcRjo1ksy7zwvcxwA%3D%3D&trackingId=XSF1GYCf6L3CCXgvzuHgqA%3D%3D&trk=flagship3_search_srp_jobs'

# Row 6
$ws.Range("A6").Value = 'Analyst Programmer'
$ws.Range("B6").Value = 'Element Materials Technology'
$ws.Range("C6").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F6").Value = 'https://www.linkedin.com/jobs/view/2596932004/?eBP=CwEAAAF6sFnSRYphO-XBYBC_9I24OFOQwsoWvvtZHkfsCkcZ4zgx6xvmyFJhRJXWRBb4BrWayHpzWZD8fk9Bk6RtaGsueDDAu7Xl9fRY8H8FkDK7lLt3QpSlbd3BVlNMpzQHzpSfpCnSE255WMf6qIFHegkQju1dR7N2_HfA6-uPnrD4d3lgPgwpZAN9Y2RIV9MFpUdLORt0nVcYVfmzNvla7Rf4xrNoevE8Jc9mBYjCoLJzwHWsIUFjm0SQhl47hf60rgALmhdYAXeNRUIQ1ZR-DelSejNWOtg4RCFqXh405HNGYUJDO5osib6XI-D0qpdwzI_iHuYYMWCqJnccznX3fKSbFAdTd7ASPhf8nOc94mT_GQVx4IlD7ANgvKD633c&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=IQ08OXcRjo1ksy7zwvcxwA%3D%3D&trackingId=8EIaaLeSqWGMp8w%2FFKfzrw%3D%3D&trk=flagship3_search_srp_jobs'

# Row 7
$ws.Range("A7").Value = 'Senior Software Engineer - Backend'
$ws.Range("B7").Value = 'HungerStation | هنقرستيشن'
$ws.Range("C7").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F7").Value = 'https://www.linkedin.com/jobs/view/2632660011/?eBP=CwEAAAF6sFnSRTVNVugTMdlPYgBqFZ4uxLfuC5wunFkmCfr-yeWCoBqrc3CRSPtrdD6iTqVXAjQ_5HOnbKGoqRBStqCr7cF-2DrqOF1G41q6cWrcwOPksSeNIIkFBxRvQc4qCPHANqWzFXVjArSkCqbW-ZraQbnGQd68nEnbMBd0gfleUX_ERrAiN3NxvA8LGlgsJmX4iBKypn8D_0biYrLvavlfUM6IuMVjVsWs16Vb0p2vI59U1dkQGb8Dk1EYpF-vofuVF_b-3--PYqoMFmjIFZJoiZBjVAC6QqnyzIo47rBP58F_XBIjfojcodv6ccnnFu9PoLcLe9M1c5EyY2B7khHb4ihoNIPyyk4hD035R1Pksdr0n3dQMIC4NQ&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=IQ08OXcRjo1ksy7zwvcxwA%3D%3D&trackingId=maVlzUTwp%2F5D1nwlOM%2BvbQ%3D%3D&trk=flagship3_search_srp_jobs'

# Row 26
$ws.Range("A26").Value = 'Kotlin Developer – Full-time, Remote'
$ws.Range("B26").Value = 'Toptal'
$ws.Range("C26").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F26").Value = 'https://www.linkedin.com/jobs/view/2648921423/?eBP=CwEAAAF6sGCRLNffWS3pcIqKajJonH2AaqFyPzJ4TpLMPT5BYqsILnNN47bGAFhDP6KQkuKyQLjgIVjWGN5nMIFKIk911wh2yOWNkZrb6aSN0bFLm75DS_DUzj7FoAa8OtVTVqJ1W1rSvsFnwGLZJtwCS9SthGwCneSlpwkuHBjD7DgYNeSx-RRoQ7TD8_YZQvog1wrjTqyywTIhuN0P3feo6RtNmU-KEkkiJLx78CsP9GPgh1CqOQYSea0cqv8NcAfaySzSc1vNT3Vl-j0mdAg8O_ze9Acl408KT9iRASHIZL9RbxJF8t3bzJ4Bu6XYQ98l8W6aLI-kRxLzi3l_Zl6PmtBYJi_0M9fO6v1lEARrOoznZZI9oRyXEH--Yh1Th4pAnQ&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=GZ2xbL4fakYgxg1w%2FA5t%2FA%3D%3D&trackingId=9wVIa3Cfw%2Fm9r4IFU7x4Pg%3D%3D&trk=flagship3_search_srp_jobs'

# Row 27
$ws.Range("A27").Value = 'Senior C++ Developer – Full-time, Remote'
$ws.Range("B27").Value = 'Toptal'
$ws.Range("C27").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F27").Value = 'https://www.linkedin.com/jobs/view/2637116275/?eBP=CwEAAAF6sGCRLL0gJcxrUp35wUAclZE6LQAeeUXI7irc05gI2ae7TCDW27l5G77cq92QpTIH0S2rqITGilVUgDdXhT8Vp38OSVtcGJGgiVAVL-0NdglDIskOjyqNXJzdjuwiSUiOBAaCNNrn5wog0Z1hwBQgNGZ3VlJhmmKU75DB6kOT6LoBSfml9RrHAxowBmduAPL_Jggr_rlwCBvFTg3V2LgTAB3hbyu6UjfEljBoW8N9q18FNbTS497CFcv1vyUO4S_0k0Qx4yx72AooU08Z3yAXCNLeszZU-SJiRTY8Dh15cPVp7xoADJsHv_YbSNTlfO68uq12pe8xhSQg5D5HnGKA85eNXhEKERx_O4n3uHVlXV8BnrzAFP0MzKN5N66yeg&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=GZ2xbL4fakYgxg1w%2FA5t%2FA%3D%3D&trackingId=KB2WDJ8XJHOGAFjH8sHjJg%3D%3D&trk=flagship3_search_srp_jobs'

# Row 28
$ws.Range("A28").Value = 'Senior Software Engineer - iOS'
$ws.Range("B28").Value = 'HungerStation | هنقرستيشن'
$ws.Range("C28").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F28").Value = 'https://www.linkedin.com/jobs/view/2631668884/?eBP=CwEAAAF6sGCRLItbkls_7TwLGQ6rsgz4PyunNjv5HGQNN-e7dOaCzWfAC4mUx-ZVpndgELD_ujCkaYSSwUKSNNqmcavVZjdGIVM99gaaCpw4SdSLRmq30YthZbe3pIYO8O2KaIsfkbQO6BfUmkHpM_1zZ1v0JeQ-FydYBizn6qXJVXY8QX1d_-JBah_S85TSHRxfeEdwxwbm9eWkvfwxF21lVogbw-Zhvjj8DT0cNX0_v_j1bTk-fI2gmdI2PxQgSGFFpjwwQJXKTKYsdRrFgLuwtDUCQTxrYR8vHIWZypaYhgk90gavBmEWGSZI3rTC17CHTAI1VPjNRJR9linC7PJPtoi_lC_XT18x7vUovB2ULZ7_wS4tKRWvZxsRL1M1YlU&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=GZ2xbL4fakYgxg1w%2FA5t%2FA%3D%3D&trackingId=IIwtSF%2BEskNLZGcx%2BdEL%2FA%3D%3D&trk=flagship3_search_srp_jobs'

# Row 29
$ws.Range("A29").Value = 'MuleSoft Developer (On-site)'
$ws.Range("B29").Value = 'Halian'
$ws.Range("C29").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("E29").Value = '
              2 days ago
 '
$ws.Range("F29").Value = 'https://www.linkedin.com/jobs/view/2643770873/?eBP=JOB_SEARCH_ORGANIC&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=GZ2xbL4fakYgxg1w%2FA5t%2FA%3D%3D&trackingId=xT3lv4vxWVcffR7nvrvHUQ%3D%3D&trk=flagship3_search_srp_jobs'

# Row 30
$ws.Range("A30").Value = 'Senior Software Engineer I - Android'
$ws.Range("B30").Value = 'talabat'
$ws.Range("C30").Value = 'Business Bay, Dubai, United Arab Emirates'
$ws.Range("F30").Value = 'https://www.linkedin.com/jobs/view/2636396773/?eBP=CwEAAAF6sGCRLFstEgneMknLFNynKuv4n2hq7Vf2FLrvwj4d4TXId0bSRyJIXaIn46GxH6Vl5cgHPk87KVrASGkOtx6E5MbldWhDGa01CJn7ev8UIj7KdB-DO3Ve_BOyNDWAmo_hmhjIG9sFZTF0uu2JWzcRSMLi2eAcPikSRapC4mVZwwcO01flEmda3hCTXSJPCotZNVYTiQi8wldIonMonQVbBPSevgwMoT2zpAKA_K_OTyF_PJ4T0t_DbiuwFS74icvDHHxYTFHFevjKAULsdniAsxHekuv3KCa2gLaVVk2KjSAdFhq60tWMHKqTwkujpz1Id_GnHITLHBdoVUXETgTsjVO-w7xg6UKOgzjdSDyAI_N1FQMhGrjT1lh2JYtC-w&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=GZ2xbL4fakYgxg1w%2FA5t%2FA%3D%3D&trackingId=b7b0Yd0sShh6bogMbXCfGQ%3D%3D&trk=flagship3_search_srp_jobs'

# Row 31
$ws.Range("A31").Value = 'Senior Software Engineer - Backend'
$ws.Range("B31").Value = 'HungerStation | هنقرستيشن'
$ws.Range("C31").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F31").Value = 'https://www.linkedin.com/jobs/view/2632660011/?eBP=CwEAAAF6sGCRLA0oWWdl2n4MX2NIapZz6HvKN7x_RkwRoYRndueWy7Fpi-HCfE2nZJDG4ed1JnMttGAaC5A7raLdhcXzyuQ550m6de2CQSrZOIMpGalrSqvWAXBifBxQoPQEcdZR38TsPnJ57hhDBnvWDWELnci6toYX78Afh3CbhFx1ufc04VYh_wojACHR8luKraRkGjaxGHgC_WctNfPIqyneW2xH4-mGSR2tIx96G2f-cCGHkcj64wzggHXwdwTz-shYCfOEPSHEihXU3-yLlNvl8NQZM6DlsxOg1Eh52dTWQgWUU-LCWCh9RNwBqAaoo-GzWJXHH9rMHApr9TT2zM03J_xFrvnLg_Z-aM6dZ2vILD3TuygIjZ1FeCQbOHAySQ&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=GZ2xbL4fakYgxg1w%2FA5t%2FA%3D%3D&trackingId=oxbwDZMI9epbV0gicXUa%2Bg%3D%3D&trk=flagship3_search_srp_jobs'

# Row 32
$ws.Range("A32").Value = 'Analyst Programmer'
$ws.Range("B32").Value = 'Element Materials Technology'
$ws.Range("C32").Value = 'Dubai, Dubai, United Arab Emirates'
$ws.Range("F32").Value = 'https://www.linkedin.com/jobs/view/2596932004/?eBP=CwEAAAF6sGCRLHmy1_cE9GM9WsDXd5AFzaMuuiLawfDlvN21zzvRRfgZWwDGOnCGSnmj8M2zakmcUqO9Qbex5WyU84XAnptTKwJERgFUKjnHlFRAFMQjxC6t3qshN3XhYhXh1vKJKhCn2aTYqn0BTyUodItcUhjKe4yVD5XQqJJ8uuwoKFisO9jb5sjc2dditp1y5sBDzWrUlrWYs9MqlUHOEKcd2C9wyal1g6rb2U9yxPxNL0RNyMIqWC1yAjPpvxU1g9HlQsWpPzCo1ikD4bKJ3pHSTtJQFrbe9hXAjedj-ED8wYLi-mmRGTY1aKBV9pEbz1UJGEY0Sv5Xgsq-huMm656_O7_8BIBww16J9oXSyAh5wNH-d3S7p7gpfgTNjpyPOOfB&recommendedFlavor=ACTIVELY_HIRING_COMPANY&refId=GZ2xbL4fakYgxg1w%2FA5t%2FA%3D%3D&trackingId=YN3cKGE2ZZjFeMSzAdIhKA%3D%3D&trk=flagship3_search_srp_jobs'

# --- Trailing near-empty rows (49 blank, 50 materialised-but-empty) to extend the used range to F50 ---
$ws.Range("A50:F50").Font.Bold = $false

# --- Selection / active-sheet bookkeeping ---
$ws.Cells.Select()
